$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(2, 8707.0, 184840.0),
    @(3, 19413.0, 377822.0),
    @(4, 30922.0, 593890.0),
    @(5, 42819.0, 708475.0),
    @(6, 55233.0, 896078.0),
    @(7, 67832.0, 1076231.0),
    @(8, 80678.0, 1268119.0),
    @(9, 93647.0, 1467730.0),
    @(10, 106969.0, 1646754.0),
    @(11, 120450.0, 1844116.0),
    @(12, 134021.0, 2051031.0),
    @(13, 147684.0, 2233631.0),
    @(14, 161472.0, 2436466.0),
    @(15, 175369.0, 2671051.0),
    @(16, 189265.0, 2858577.0),
    @(17, 203275.0, 3046804.0),
    @(18, 217519.0, 3264744.0),
    @(19, 231958.0, 3460276.0),
    @(20, 246402.0, 3644073.0),
    @(21, 260908.0, 3857788.0),
    @(22, 275456.0, 4067654.0),
    @(23, 290040.0, 4277734.0),
    @(24, 304664.0, 4423179.0),
    @(25, 319333.0, 4702565.0),
    @(26, 334127.0, 4913304.0),
    @(27, 348969.0, 5097964.0),
    @(28, 363833.0, 5327853.0),
    @(29, 378734.0, 5535992.0),
    @(30, 393630.0, 5716145.0),
    @(31, 408611.0, 5940235.0),
    @(32, 423585.0, 6185374.0),
    @(33, 438549.0, 6349080.0),
    @(34, 453671.0, 6581252.0),
    @(35, 469056.0, 6786757.0),
    @(36, 484477.0, 7000840.0),
    @(37, 499912.0, 7179308.0),
    @(38, 515349.0, 7359726.0),
    @(39, 530817.0, 7627950.0),
    @(40, 546303.0, 7836747.0),
    @(41, 561767.0, 8028772.0),
    @(42, 577303.0, 8225347.0),
    @(43, 592915.0, 8428901.0),
    @(44, 608509.0, 8662058.0),
    @(45, 624110.0, 8779038.0),
    @(46, 639717.0, 8888790.0),
    @(47, 655374.0, 9120570.0),
    @(48, 671028.0, 9350681.0),
    @(49, 686669.0, 9562216.0),
    @(50, 702366.0, 9733867.0),
    @(51, 718183.0, 9987405.0),
    @(52, 734036.0, 10207740.0),
    @(53, 749901.0, 10421139.0),
    @(54, 765793.0, 10622564.0),
    @(55, 781642.0, 10820900.0),
    @(56, 797563.0, 11092024.0),
    @(57, 813485.0, 11286119.0),
    @(58, 829336.0, 11459497.0),
    @(59, 845263.0, 11728166.0),
    @(60, 861219.0, 11919840.0),
    @(61, 877198.0, 12158043.0),
    @(62, 893180.0, 12317249.0),
    @(63, 909116.0, 12525388.0),
    @(64, 925113.0, 12794331.0),
    @(65, 941152.0, 12984585.0),
    @(66, 957107.0, 13166542.0),
    @(67, 973317.0, 13412152.0),
    @(68, 989741.0, 13649089.0),
    @(69, 1006119.0, 13870665.0),
    @(70, 1022533.0, 14094994.0),
    @(71, 1038931.0, 14284949.0),
    @(72, 1055397.0, 14476478.0),
    @(73, 1071807.0, 14773405.0),
    @(74, 1088228.0, 15017852.0),
    @(75, 1104713.0, 15138125.0),
    @(76, 1121159.0, 15429237.0),
    @(77, 1137659.0, 15602436.0),
    @(78, 1154067.0, 15817315.0),
    @(79, 1170604.0, 15936065.0),
    @(80, 1187101.0, 16204153.0),
    @(81, 1203605.0, 16481357.0),
    @(82, 1220073.0, 16656207.0),
    @(83, 1236620.0, 16834513.0),
    @(84, 1253193.0, 17111025.0),
    @(85, 1269780.0, 17383593.0),
    @(86, 1286456.0, 17619085.0),
    @(87, 1303004.0, 17809305.0),
    @(88, 1319620.0, 18029358.0),
    @(89, 1336218.0, 18342665.0),
    @(90, 1352864.0, 18497433.0),
    @(91, 1369474.0, 18703168.0),
    @(92, 1386114.0, 18959050.0),
    @(93, 1402774.0, 19169429.0),
    @(94, 1419337.0, 19476612.0),
    @(95, 1436097.0, 19605011.0),
    @(96, 1452715.0, 19624777.0),
    @(97, 1469407.0, 19837115.0),
    @(98, 1486085.0, 20394364.0),
    @(99, 1502727.0, 20586474.0),
    @(100, 1519480.0, 20806416.0),
    @(101, 1536378.0, 21062734.0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
